# Slide 4 ("object 19" textbox, shape #14 in z-order): the office-hours
# entry for nada.mohyudeen@guc.edu.eg currently reads
#   "Monday 3rd  or by appointment  Saturday"
# and needs to read
#   "Sunday 3rd  or by appointment  Saturday"
# The source run "Monday 3" is shortened to "3" and a brand-new run
# "Sunday " is inserted immediately before it (mirrors a user retyping the
# day name in the PowerPoint UI: insert the new word, then delete the old
# one), leaving the "rd" (superscript), " or by appointment ", "Saturday"
# runs untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(14)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Paragraph 3 is "Monday 3rd  or by appointment  Saturday"; its first run
# is "Monday 3" (the rest - "rd", "  or", " ", "by", " ", "appointment ",
# " ", "Saturday" - are later runs we must not disturb).
$para = $tr.Paragraphs(3, 1)
$run = $para.Runs(1, 1)

# Shrink "Monday 3" down to just "3", keeping its original run/formatting
# (no smtClean attribute, same as the source run it came from).
$run.Text = "3"

# Re-acquire the (now 1-character) run and type the new word in front of
# it as a separate run - this is exactly what typing fresh text in the
# PowerPoint editor produces (a distinct run carrying dirty/smtClean).
$run = $tr.Paragraphs(3, 1).Runs(1, 1)
[void]$run.InsertBefore("Sunday ")

# The textbox uses <a:spAutoFit/>, so PowerPoint relays it out and grows
# it by a touch when the wrapped text reflows with the new word (width
# 1990089 EMU is unchanged; only the autofit height grows to 698781 EMU).
$shp.Height = 55.02212717634487
